$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 239
$ws.Range("J6").Value = 299
$ws.Range("L6").Value = 897
$ws.Range("N6").Value = -1121
$ws.Range("H64").Value = 6861.086
$ws.Range("J64").Value = 8285.571
$ws.Range("L64").Value = 8285.571
$ws.Range("N64").Value = -8781.571
$ws.Range("H67").Value = 6861.086
$ws.Range("J67").Value = 8285.571
$ws.Range("L67").Value = 8285.571
$ws.Range("N67").Value = -10001.571
$ws.Range("H87").Value = 41449.5
$ws.Range("J87").Value = 41449.5
$ws.Range("L87").Value = 41449.5
$ws.Range("N87").Value = -43945.5
$ws.Range("H90").Value = 41449.5
$ws.Range("J90").Value = 41449.5
$ws.Range("L90").Value = 124348.5
$ws.Range("N90").Value = -136828.5
$ws.Range("H137").Value = 5736.724
$ws.Range("I137").Value = 3607
$ws.Range("K137").Value = 10821
$ws.Range("M137").Value = -8271
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4552717
$ws.Range("I32").Value = 4883269.5
$ws.Range("K32").Value = 4883269.5
$ws.Range("M32").Value = -4882982.5
$ws.Range("H34").Value = 72499.56
$ws.Range("I34").Value = 84284.86
$ws.Range("K34").Value = 84284.86
$ws.Range("M34").Value = -84013.86
$ws.Range("H62").Value = 41400
$ws.Range("J62").Value = 41400
$ws.Range("L62").Value = 41400
$ws.Range("N62").Value = -42648
$ws.Range("H65").Value = 41400
$ws.Range("J65").Value = 41400
$ws.Range("L65").Value = 124200
$ws.Range("N65").Value = -130440
$ws.Range("H102").Value = 26319762
$ws.Range("I102").Value = 50003490
$ws.Range("J102").Value = 4511
$ws.Range("K102").Value = 50003490
$ws.Range("L102").Value = 4511
$ws.Range("M102").Value = -50001868
$ws.Range("N102").Value = -7755
$ws.Range("H122").Value = 4631.476
$ws.Range("I122").Value = 2098.25
$ws.Range("K122").Value = 6294.75
$ws.Range("M122").Value = -3844.75
$ws.Range("H123").Value = 50233.5
$ws.Range("J123").Value = 50567
$ws.Range("L123").Value = 50567
$ws.Range("N123").Value = -60367
$ws.Range("H126").Value = 4399.1665
$ws.Range("I126").Value = 4399.1665
$ws.Range("K126").Value = 13197.4995
$ws.Range("M126").Value = -10727.4995
$ws.Range("H132").Value = 6999.8164
$ws.Range("I132").Value = 5994.143
$ws.Range("J132").Value = 8340.714
$ws.Range("K132").Value = 17982.429
$ws.Range("L132").Value = 25022.142
$ws.Range("M132").Value = -15452.429
$ws.Range("N132").Value = -30082.142
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 27683.21
$ws.Range("I86").Value = 37701.48
$ws.Range("K86").Value = 37701.48
$ws.Range("M86").Value = -36578.48
$ws.Range("H89").Value = 27683.21
$ws.Range("I89").Value = 37701.48
$ws.Range("K89").Value = 188507.4
$ws.Range("M89").Value = -182891.4
$ws.Range("H105").Value = 5916.59
$ws.Range("J105").Value = 5143.643
$ws.Range("L105").Value = 5143.643
$ws.Range("N105").Value = -8637.643
$ws.Range("H122").Value = 70233.5
$ws.Range("J122").Value = 70233.5
$ws.Range("L122").Value = 70233.5
$ws.Range("N122").Value = -80033.5
$ws.Range("H128").Value = 3908.875
$ws.Range("I128").Value = 3908.875
$ws.Range("K128").Value = 11726.625
$ws.Range("M128").Value = -9236.625
$ws.Range("H134").Value = 4977.0225
$ws.Range("I134").Value = 1809.0333
$ws.Range("K134").Value = 5427.0999
$ws.Range("M134").Value = -2892.0999
$ws.Range("H139").Value = 67165.5
$ws.Range("I139").Value = 64997
$ws.Range("J139").Value = 67599.2
$ws.Range("K139").Value = 64997
$ws.Range("L139").Value = 67599.2
$ws.Range("M139").Value = -59857
$ws.Range("N139").Value = -77879.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 840.5238000000001
$ws.Range("I5").Value = 269.66666
$ws.Range("J5").Value = 1601.6666
$ws.Range("K5").Value = 269.66666
$ws.Range("L5").Value = 1601.6666
$ws.Range("M5").Value = -157.66666
$ws.Range("N5").Value = -1825.6666
$ws.Range("H31").Value = 7467.425
$ws.Range("I31").Value = 2003.7059
$ws.Range("K31").Value = 2003.7059
$ws.Range("M31").Value = -1708.7059
$ws.Range("H34").Value = 7467.425
$ws.Range("I34").Value = 2003.7059
$ws.Range("K34").Value = 2003.7059
$ws.Range("M34").Value = -1801.7059
$ws.Range("H53").Value = 54958.25
$ws.Range("J53").Value = 54958.25
$ws.Range("L53").Value = 54958.25
$ws.Range("N53").Value = -56172.25
$ws.Range("H100").Value = 44768
$ws.Range("J100").Value = 49536
$ws.Range("L100").Value = 49536
$ws.Range("N100").Value = -51700
$ws.Range("H134").Value = 5133.718
$ws.Range("I134").Value = 1540.9
$ws.Range("J134").Value = 8915.632
$ws.Range("K134").Value = 4622.700000000001
$ws.Range("L134").Value = 26746.896
$ws.Range("M134").Value = -2087.700000000001
$ws.Range("N134").Value = -31816.896
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27770270
$ws.Range("I4").Value = 39326156
$ws.Range("K4").Value = 117978468
$ws.Range("M4").Value = -117978356
$ws.Range("H5").Value = 2274.1538
$ws.Range("I5").Value = 1173.8889
$ws.Range("K5").Value = 3521.6667
$ws.Range("M5").Value = -3409.6667
$ws.Range("H119").Value = 2607.25
$ws.Range("I119").Value = 214.5
$ws.Range("J119").Value = 5000
$ws.Range("K119").Value = 643.5
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = 4194.5
$ws.Range("N119").Value = -24676
$ws.Range("H134").Value = 48383.957
$ws.Range("I134").Value = 59046.168
$ws.Range("K134").Value = 177138.504
$ws.Range("M134").Value = -172068.504
$ws.Range("H135").Value = 2274.1538
$ws.Range("I135").Value = 1173.8889
$ws.Range("K135").Value = 10565.0001
$ws.Range("M135").Value = -8030.000099999999
$ws.Range("H138").Value = 78783.14
$ws.Range("I138").Value = 89663.664
$ws.Range("K138").Value = 268990.992
$ws.Range("M138").Value = -263850.992
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2800
$ws.Range("I80").Value = 2800
$ws.Range("K80").Value = 2800
$ws.Range("M80").Value = -1802
$ws.Range("H83").Value = 2800
$ws.Range("I83").Value = 2800
$ws.Range("K83").Value = 14000
$ws.Range("M83").Value = -9008
$ws.Range("H102").Value = 4161.8276
$ws.Range("I102").Value = 3871.7693
$ws.Range("J102").Value = 6675.6665
$ws.Range("K102").Value = 3871.7693
$ws.Range("L102").Value = 6675.6665
$ws.Range("M102").Value = -2249.7693
$ws.Range("N102").Value = -9919.666499999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2333333.2
$ws.Range("I20").Value = 2250000
$ws.Range("K20").Value = 2250000
$ws.Range("M20").Value = -2249774
$ws.Range("H74").Value = 20000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 20000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 6567.5815
$ws.Range("I132").Value = 3152.9473
$ws.Range("J132").Value = 9270.833000000001
$ws.Range("K132").Value = 9458.841899999999
$ws.Range("L132").Value = 27812.499
$ws.Range("M132").Value = -6928.841899999999
$ws.Range("N132").Value = -32872.499
$ws.Range("H136").Value = 11991.718
$ws.Range("J136").Value = 19523.8
$ws.Range("L136").Value = 58571.39999999999
$ws.Range("N136").Value = -63671.39999999999
$ws.Range("H139").Value = 78672.73
$ws.Range("J139").Value = 82040
$ws.Range("L139").Value = 82040
$ws.Range("N139").Value = -92320
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166671460
$ws.Range("I62").Value = 142862610
$ws.Range("J62").Value = 333333340
$ws.Range("K62").Value = 142862610
$ws.Range("L62").Value = 333333340
$ws.Range("M62").Value = -142861986
$ws.Range("N62").Value = -333334588
$ws.Range("H65").Value = 166671460
$ws.Range("I65").Value = 142862610
$ws.Range("J65").Value = 333333340
$ws.Range("K65").Value = 714313050
$ws.Range("L65").Value = 1666666700
$ws.Range("M65").Value = -714309930
$ws.Range("N65").Value = -1666672940
$ws.Range("H107").Value = 12346713
$ws.Range("I107").Value = 531.1429000000001
$ws.Range("J107").Value = 25642602
$ws.Range("K107").Value = 1593.4287
$ws.Range("L107").Value = 76927806
$ws.Range("M107").Value = 326.5712999999998
$ws.Range("N107").Value = -76931646
$ws.Range("H126").Value = 1273.75
$ws.Range("I126").Value = 1284.2858
$ws.Range("K126").Value = 3852.8574
$ws.Range("M126").Value = -1382.8574
$ws.Range("H129").Value = 99999
$ws.Range("J129").Value = 99999
$ws.Range("L129").Value = 99999
$ws.Range("N129").Value = -109999
$ws.Range("H132").Value = 5871.2334
$ws.Range("I132").Value = 8810.462
$ws.Range("J132").Value = 3623.5881
$ws.Range("K132").Value = 26431.386
$ws.Range("L132").Value = 10870.7643
$ws.Range("M132").Value = -23901.386
$ws.Range("N132").Value = -15930.7643
$ws.Range("H136").Value = 4046.2856
$ws.Range("I136").Value = 1504.4286
$ws.Range("J136").Value = 6588.143
$ws.Range("K136").Value = 4513.2858
$ws.Range("L136").Value = 19764.429
$ws.Range("M136").Value = -1963.2858
$ws.Range("N136").Value = -24864.429
